$wb = $excel.ActiveWorkbook

# Add the data row on the Data sheet
$data = $wb.Worksheets.Item("Data")
$data.Range("A2").Value = "bilateral-unspecified"
$data.Range("B2").Value = "Bilateral, unspecified"
$data.Range("C2").Value = 2015
$data.Range("D2").Value = 125175540

# Fix the "Units of measure" note on the Notes sheet
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"
